# DOMA-10686: Fix export drops
# The contacts export template had 10 rows pre-provisioned for data
# (header + 2 sample rows + 7 blank placeholder rows), which made Excel
# pad the export with extra empty rows. Trim the sheet back down to the
# 3 rows that are actually used (header + the two sample/template rows),
# and let the last header cell (H3) wrap its text instead of overflowing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the cursor where the author last left it (inside the area that is
# about to be removed) before deleting the now-unused trailing rows.
[void]$ws.Range("F8").Select()

# Drop the 7 blank placeholder rows (4-10); only the header + two sample
# rows are needed.
$ws.Rows("4:10").Delete()

# Let the last cell in the remaining header row wrap its text.
$ws.Range("H3").WrapText = $true

# Drop the page margins so the exported sheet prints edge-to-edge.
$ws.PageSetup.LeftMargin = 0
$ws.PageSetup.RightMargin = 0
$ws.PageSetup.TopMargin = 0
$ws.PageSetup.BottomMargin = 0
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
